# g20.1 e g20.2 - programação da nova fonte
# Replace the IDHM (2022) table with the new-source (2010) table, re-sorted,
# and give the header row a thin box border + top vertical alignment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New dataset (replaces the old 2022 figures with the 2010 ones, re-sorted
# ascending by value, "Acre" dropped and "Piauí" added).
$data = @(
  @("Sergipe",  "IDHM", "31/12/2010", 0.665, "20º"),
  @("Bahia",    "IDHM", "31/12/2010", 0.66,  "22º"),
  @("Paraíba",  "IDHM", "31/12/2010", 0.658, "23º"),
  @("Pará",     "IDHM", "31/12/2010", 0.646, "24º"),
  @("Piauí",    "IDHM", "31/12/2010", 0.646, "25º"),
  @("Maranhão", "IDHM", "31/12/2010", 0.639, "26º"),
  @("Alagoas",  "IDHM", "31/12/2010", 0.631, "27º"),
  @("Brasil",   "IDHM", "31/12/2010", 0.727, ""),
  @("Nordeste", "IDHM", "31/12/2010", 0.667, "")
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $data[$i][0]
  $ws.Cells.Item($row, 2).Value = $data[$i][1]
  $ws.Cells.Item($row, 3).Value = $data[$i][2]
  $ws.Cells.Item($row, 4).Value = $data[$i][3]
  $ws.Cells.Item($row, 5).Value = $data[$i][4]
}

# Header row: add a thin box border around every cell and switch the
# vertical alignment to top (horizontal center is already set).
$header = $ws.Range("A1:E1")
$header.Borders.LineStyle = 1
$header.VerticalAlignment = -4160

# Restore Excel's default page margins (0.75in/1in/0.5in).
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36
